$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores these numeric-looking stats as text (numberStoredAsText),
# so write each updated value straight into the cells it changes in. Setting
# NumberFormat to "@" immediately before the write keeps the new value typed
# as text (matching the existing convention) without touching any cell whose
# value is not actually changing.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue "C2" "17"
Set-TextValue "D2" "11"
Set-TextValue "E2" "3"

Set-TextValue "C3" "13"
Set-TextValue "D3" "14"
Set-TextValue "E3" "1"

Set-TextValue "C5" "0"
Set-TextValue "D5" "1"
Set-TextValue "E5" "0"

Set-TextValue "C6" "5"
Set-TextValue "E6" "0"

Set-TextValue "C7" "21"
Set-TextValue "D7" "18"
Set-TextValue "F7" "0"

Set-TextValue "C8" "10"
Set-TextValue "D8" "10"
Set-TextValue "E8" "0"
Set-TextValue "F8" "1"

Set-TextValue "C9" "30"
Set-TextValue "D9" "27"
Set-TextValue "E9" "2"

Set-TextValue "C10" "10"
Set-TextValue "D10" "6"
Set-TextValue "E10" "1"
